$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (month labels)
$ws.Range("B1").Value = "Series Index Apr"
$ws.Range("C1").Value = "Series Index Mar"

# Row 2 - Manufacturing PMI
$ws.Range("B2").Value = "'49.2"
$ws.Range("C2").Value = "'50.3"
$ws.Range("D2").Value = "'-1.1"
$ws.Range("F2").Value = "From Growing"
$ws.Range("G2").Value = "'1"

# Row 3 - New Orders
$ws.Range("B3").Value = "'49.1"
$ws.Range("C3").Value = "'51.4"
$ws.Range("D3").Value = "'-2.3"
$ws.Range("F3").Value = "From Growing"
$ws.Range("G3").Value = "'1"

# Row 4 - Production
$ws.Range("B4").Value = "'51.3"
$ws.Range("C4").Value = "'54.6"
$ws.Range("D4").Value = "'-3.3"

# Row 5 - Employment
$ws.Range("B5").Value = "'48.6"
$ws.Range("C5").Value = "'47.4"
$ws.Range("D5").Value = "'+1.2"
$ws.Range("F5").Value = "Slower"
$ws.Range("G5").Value = "'7"

# Row 6 - Supplier Deliveries
$ws.Range("B6").Value = "'48.9"
$ws.Range("C6").Value = "'49.9"
$ws.Range("D6").Value = "'-1.0"
$ws.Range("F6").Value = "Faster"
$ws.Range("G6").Value = "'2"

# Row 7 - Inventories
$ws.Range("B7").Value = "'48.2"
$ws.Range("C7").Value = "'48.2"
$ws.Range("D7").Value = "'0.0"
$ws.Range("F7").Value = "Same"
$ws.Range("G7").Value = "'15"

# Row 8 - Customers' Inventories
$ws.Range("B8").Value = "'47.8"
$ws.Range("C8").Value = "'44.0"
$ws.Range("D8").Value = "'+3.8"

# Row 9 - Prices
$ws.Range("B9").Value = "'60.9"
$ws.Range("C9").Value = "'55.8"
$ws.Range("D9").Value = "'+5.1"
$ws.Range("E9").Value = "Increasing"
$ws.Range("F9").Value = "Faster"
$ws.Range("G9").Value = "'4"

# Row 10 - Backlog of Orders
$ws.Range("B10").Value = "'45.4"
$ws.Range("C10").Value = "'46.3"
$ws.Range("D10").Value = "'-0.9"
$ws.Range("G10").Value = "'19"

# Row 11 - New Export Orders
$ws.Range("B11").Value = "'48.7"
$ws.Range("C11").Value = "'51.6"
$ws.Range("D11").Value = "'-2.9"
$ws.Range("F11").Value = "From Growing"
$ws.Range("G11").Value = "'1"

# Row 12 - Imports
$ws.Range("B12").Value = "'51.9"
$ws.Range("C12").Value = "'53.0"
$ws.Range("D12").Value = "'-1.1"
$ws.Range("E12").Value = "Growing"
$ws.Range("F12").Value = "Slower"
$ws.Range("G12").Value = "'4"

# Row 13 - OVERALL ECONOMY
$ws.Range("E13").Value = "Growing"
$ws.Range("F13").Value = "Slower"
$ws.Range("G13").Value = "'48"

# Row 14 - Manufacturing Sector
$ws.Range("F14").Value = "From Growing"
$ws.Range("G14").Value = "'1"
